$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 56: new TRADING_ATTEMPT log entry for SUI
$ws.Range("A56").Value = "2025-10-03T01:24:18.031794"
$ws.Range("B56").Value = "TRADING_ATTEMPT"
$ws.Range("C56").Value = "SUI"
$ws.Range("D56").Value = "UNKNOWN"
$ws.Range("E56").Value = 3.546780152288974
$ws.Range("K56").Value = "ATTEMPT"
$ws.Range("L56").Value = "Attempting trade 1/1"

# Row 57: resulting POSITION_FAILED log entry for SUI
$ws.Range("A57").Value = "2025-10-03T01:24:19.927512"
$ws.Range("B57").Value = "POSITION_FAILED"
$ws.Range("C57").Value = "SUI"
$ws.Range("D57").Value = "UNKNOWN"
$ws.Range("K57").Value = "FAILED"
$ws.Range("L57").Value = "Trade execution failed for trade 1"
